$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost")

# Row 13: extra note about postcard specs (same plain style as surrounding cells)
$ws.Range("B13").Value = "12 inches * 7.5 inches, 100 copies, $150, a day, $20 for putting address"

# New postcard / mailing-label cost breakdown block (rows 16-20)
$ws.Range("A17").Value = "postcard"
$ws.Range("B17").Value = 150

$ws.Range("A19").Value = "mailing label"
$ws.Range("B19").Value = 20

$ws.Range("A16").Value = "Postcard"
$ws.Range("B16").Value = 100
$ws.Range("A16:B16").Font.Bold = $true

$ws.Range("A18").Value = "postage stamp"
$ws.Range("B18").Formula = "=0.6*B16"

$ws.Range("A17:B19").Borders.LineStyle = 1

$ws.Range("A20").Value = "Total"
$ws.Range("B20").Formula = "=SUM(B17:B19)"
$ws.Range("B20").Font.Bold = $true

[void]$ws.Range("C5").Select()
